$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
# Date Solved for the new entries: 2025-09-17 (serial 45917)
$dateSolvedSerial = 45917

# Problem names (column C) entered first for the three new rows
$ws.Range("C76").Value = "Maximum Size Subarray Sum Equals K"
$ws.Range("C77").Value = "Maximum Subarray"
$ws.Range("C78").Value = "Maximum Sum Circular Subarray"

# Row 76 -> ID 75: Maximum Size Subarray Sum Equals K
$ws.Range("B76").Value = "Prefix Sum"
$ws.Range("D76").Value = "Medium"
$ws.Range("E76").Value = "Done"
$ws.Range("F76").Value = $dateSolvedSerial
$ws.Range("G76").Value = "O(n)"
$ws.Range("H76").Value = "O(n)"
$ws.Range("I76").Value = "Prefix Sum + HashMap"

# Row 77 -> ID 76: Maximum Subarray
$ws.Range("B77").Value = "Kadane Algo"
$ws.Range("D77").Value = "Medium"
$ws.Range("E77").Value = "Done"
$ws.Range("F77").Value = $dateSolvedSerial
$ws.Range("G77").Value = "O(n)"
$ws.Range("H77").Value = "O(1)"
$ws.Range("I77").Value = "Kadane's Algorithm"

# Row 78 -> ID 77: Maximum Sum Circular Subarray
$ws.Range("B78").Value = "Kadane Algo"
$ws.Range("D78").Value = "Medium"
$ws.Range("E78").Value = "Done"
$ws.Range("F78").Value = $dateSolvedSerial
$ws.Range("G78").Value = "O(n)"
$ws.Range("H78").Value = "O(1)"
$ws.Range("I78").Value = "Kadane's Algorithm"

# Apply the same "Date Solved" number format already used for F72:F75 by
# copying the formatting from the adjacent existing date cell
$ws.Range("F75").Copy()
$ws.Range("F76:F78").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Reflect the updated column B width and selection/scroll state from the diff
$ws.Columns.Item(2).ColumnWidth = 14.75
$ws.Range("C79").Select()
$ws.Application.ActiveWindow.ScrollRow = 63
